$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Update the date in the header line:
#    "Apóstoles, Misiones, 28" -> "Apóstoles, Misiones, 29"
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Apóstoles, Misiones, 28", $false, $false, $false, $false, $false, $true, 1, $false, "Apóstoles, Misiones, 29", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from right after
#    "en Sistemas de Computación" to right after the final
#    "." that ends the "... Visual Paradigm." paragraph.
#
#    Bookmarks.Add with the name "_GoBack" replaces any existing
#    bookmark of the same name, so the old one disappears
#    automatically once the new one is created.
#
#    A collapsed (zero-length) Range placed exactly at the end of a
#    paragraph's text cannot be used directly as the bookmark target,
#    so a temporary single character is inserted, wrapped with the
#    bookmark, and then removed again - leaving the now-collapsed
#    bookmark correctly anchored at that position.
# -----------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("Visual Paradigm.") | Out-Null

$insertPos = $target.End
$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter("X") | Out-Null

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$markerRange = $d.Range($insertPos, $insertPos + 1)
$markerRange.Text = ""

# -----------------------------------------------------------------
# 3. Update the footer PAGE field cached result from "3" to "7".
# -----------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.ClearFormatting()
$footer.Range.Find.Execute("3", $false, $false, $false, $false, $false, $true, 1, $false, "7", 2) | Out-Null

Write-Output "Edit complete."
